$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant columns shared by every data row in this block (A..K, Q, T).
$common = @{
    A = 11
    B = 'Vega Monumental Concepción'
    C = 'Bíobío'
    E = 8
    F = 'Fruta'
    G = 100101
    H = 'Berries'
    I = 100112025
    J = 'Frutilla'
    K = 'Sin especificar'
    Q = '$/bandeja 7 kilos'
    T = 7
}

$rows = @(
    @{Row=165; D=44511; L='Primera'; M=600; N=7500; O=8000; P=7792; R='Provincia de Melipilla'; S=1113},
    @{Row=166; D=44511; L='Segunda'; M=450; N=6000; O=6500; P=6222; R='Provincia de Melipilla'; S=889},
    @{Row=167; D=44306; L='Especial'; M=100; N=10000; O=10000; P=10000; R='Provincia de Melipilla'; S=1429},
    @{Row=168; D=44306; L='Primera'; M=100; N=8000; O=8000; P=8000; R='Provincia de Melipilla'; S=1143},
    @{Row=169; D=44211; L='Especial'; M=200; N=10000; O=10000; P=10000; R='Provincia de Melipilla'; S=1429},
    @{Row=170; D=44211; L='Primera'; M=150; N=8000; O=8000; P=8000; R='Provincia de Melipilla'; S=1143},
    @{Row=171; D=44211; L='Segunda'; M=50; N=7000; O=7000; P=7000; R='Provincia de Melipilla'; S=1000},
    @{Row=172; D=44469; L='Primera'; M=50; N=17000; O=17000; P=17000; R='Provincia de Melipilla'; S=2429},
    @{Row=173; D=44215; L='Especial'; M=150; N=12000; O=12000; P=12000; R='Provincia de Melipilla'; S=1714},
    @{Row=174; D=44215; L='Primera'; M=100; N=10000; O=10000; P=10000; R='Provincia de Melipilla'; S=1429},
    @{Row=175; D=44215; L='Segunda'; M=100; N=8000; O=8000; P=8000; R='Provincia de Melipilla'; S=1143},
    @{Row=176; D=44504; L='Especial'; M=100; N=9000; O=9000; P=9000; R='Provincia de Melipilla'; S=1286},
    @{Row=177; D=44504; L='Primera'; M=100; N=7000; O=7000; P=7000; R='Provincia de Melipilla'; S=1000},
    @{Row=178; D=44504; L='Segunda'; M=50; N=6000; O=6000; P=6000; R='Provincia de Melipilla'; S=857},
    @{Row=179; D=44505; L='Especial'; M=200; N=8000; O=8000; P=8000; R='Provincia de Melipilla'; S=1143},
    @{Row=180; D=44505; L='Primera'; M=200; N=6500; O=6500; P=6500; R='Provincia de Melipilla'; S=929},
    @{Row=181; D=44505; L='Segunda'; M=100; N=5500; O=5500; P=5500; R='Provincia de Melipilla'; S=786},
    @{Row=182; D=44168; L='Especial'; M=100; N=10000; O=10000; P=10000; R='Provincia de Melipilla'; S=1429},
    @{Row=183; D=44168; L='Primera'; M=100; N=8000; O=8000; P=8000; R='Provincia de Melipilla'; S=1143},
    @{Row=184; D=44168; L='Segunda'; M=50; N=7000; O=7000; P=7000; R='Provincia de Melipilla'; S=1000},
    @{Row=185; D=44484; L='Especial'; M=100; N=12000; O=12000; P=12000; R='Provincia de Melipilla'; S=1714},
    @{Row=186; D=44484; L='Primera'; M=100; N=9000; O=9000; P=9000; R='Provincia de Melipilla'; S=1286},
    @{Row=187; D=44328; L='Especial'; M=100; N=15000; O=15000; P=15000; R='Provincia de Melipilla'; S=2143},
    @{Row=188; D=44328; L='Primera'; M=100; N=12000; O=12000; P=12000; R='Provincia de Melipilla'; S=1714},
    @{Row=189; D=44217; L='Especial'; M=100; N=9000; O=9000; P=9000; R='Región del Maule'; S=1286},
    @{Row=190; D=44217; L='Primera'; M=100; N=8000; O=8000; P=8000; R='Región del Maule'; S=1143},
    @{Row=191; D=44217; L='Segunda'; M=100; N=7000; O=7000; P=7000; R='Región del Maule'; S=1000},
    @{Row=192; D=44509; L='Primera'; M=450; N=7500; O=8000; P=7722; R='Provincia de Melipilla'; S=1103},
    @{Row=193; D=44509; L='Segunda'; M=250; N=6000; O=6500; P=6260; R='Provincia de Melipilla'; S=894},
    @{Row=194; D=44421; L='Primera'; M=50; N=22000; O=22000; P=22000; R='Provincia de Melipilla'; S=3143},
    @{Row=195; D=44433; L='Segunda'; M=50; N=24000; O=24000; P=24000; R='Provincia de Melipilla'; S=3429},
    @{Row=196; D=44491; L='Especial'; M=270; N=9000; O=10000; P=9556; R='Provincia de Melipilla'; S=1365},
    @{Row=197; D=44491; L='Primera'; M=380; N=6500; O=7000; P=6737; R='Provincia de Melipilla'; S=962},
    @{Row=198; D=44491; L='Segunda'; M=200; N=5500; O=5500; P=5500; R='Provincia de Melipilla'; S=786}
)

# New dimension grows to row 198 once these rows are populated; rows 197-198
# did not exist before, so every column (not only the ones that vary) must be written.
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $common.A
    $ws.Cells.Item($r.Row, 2).Value = $common.B
    $ws.Cells.Item($r.Row, 3).Value = $common.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $common.E
    $ws.Cells.Item($r.Row, 6).Value = $common.F
    $ws.Cells.Item($r.Row, 7).Value = $common.G
    $ws.Cells.Item($r.Row, 8).Value = $common.H
    $ws.Cells.Item($r.Row, 9).Value = $common.I
    $ws.Cells.Item($r.Row, 10).Value = $common.J
    $ws.Cells.Item($r.Row, 11).Value = $common.K
    $ws.Cells.Item($r.Row, 12).Value = $r.L
    $ws.Cells.Item($r.Row, 13).Value = $r.M
    $ws.Cells.Item($r.Row, 14).Value = $r.N
    $ws.Cells.Item($r.Row, 15).Value = $r.O
    $ws.Cells.Item($r.Row, 16).Value = $r.P
    $ws.Cells.Item($r.Row, 17).Value = $common.Q
    $ws.Cells.Item($r.Row, 18).Value = $r.R
    $ws.Cells.Item($r.Row, 19).Value = $r.S
    $ws.Cells.Item($r.Row, 20).Value = $common.T
}

# D column (Fecha) keeps the date number format used by the rest of the column.
$ws.Range("D165:D198").NumberFormat = $ws.Range("D164").NumberFormat
